$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to the new "custom accuracy" (2-decimal-place) readings.
$ws.Range("B5").Value = 3.9
$ws.Range("C5").Value = 2.63
$ws.Range("D5").Value = 0.01
$ws.Range("E5").Value = 6.64
$ws.Range("F5").Value = 5.86
$ws.Range("G5").Value = 2.93
$ws.Range("H5").Value = 15.78
$ws.Range("I5").Value = 3.69
$ws.Range("J5").Value = 1.88
$ws.Range("K5").Value = 3.13
$ws.Range("L5").Value = 2.79
$ws.Range("M5").Value = 2.72
$ws.Range("N5").Value = 0.86
$ws.Range("O5").Value = 2.39
$ws.Range("P5").Value = 4.12
$ws.Range("Q5").Value = 2.03
$ws.Range("R5").Value = 0.26
$ws.Range("S5").Value = 0.01
$ws.Range("T5").Value = 32.61
$ws.Range("U5").Value = 7.65
$ws.Range("V5").Value = 2.6
$ws.Range("W5").Value = 5.31
$ws.Range("X5").Value = 2.5
$ws.Range("Y5").Value = 0.35
$ws.Range("Z5").Value = 7.2
$ws.Range("AA5").Value = 2.12
$ws.Range("AB5").Value = 1.71
$ws.Range("AC5").Value = 2.06
$ws.Range("AD5").Value = 3.4
$ws.Range("AE5").Value = 0.52
$ws.Range("AF5").Value = 14.42
$ws.Range("AG5").Value = 1.23
$ws.Range("AH5").Value = 2.87

# Row 6 (the last data row) is removed entirely, shrinking the used range
# from A1:AH6 down to A1:AH5.
$ws.Rows(6).Delete()
